$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Target widths from diff: col A = 24.16796875, col B = 95.5078125 characters.
# This runtime rounds ColumnWidth to the nearest 1/7 of a character (pixel-grid snapping),
# so we pick the ColumnWidth value that lands closest to the desired stored width.
$ws.Columns.Item(1).ColumnWidth = 23.428571428571427
$ws.Columns.Item(2).ColumnWidth = 94.85714285714286

# --- Insert two new rows into the table (140 rows -> 142 rows) ---
# New "GDP" entry inserted before the old row 87 ("gender pay gap"), shifting rows 87-140 down by one.
$ws.Rows.Item(87).Insert()
# New "Mikrometer" (µm) entry inserted before the old (already shifted) row 140 ("µg/m³"), shifting further.
$ws.Rows.Item(140).Insert()

# Newly inserted rows have no cell style; copy the standard data-row format (borders/fill/font)
# from a neighbouring data row so they match the rest of the table.
$ws.Range("A86:C86").Copy()
$ws.Range("A87:C87").PasteSpecial(-4122)
$ws.Range("A139:C139").Copy()
$ws.Range("A140:C140").PasteSpecial(-4122)

# --- Append three brand-new rows at the end (143-145) and give them the same data-row format ---
$ws.Range("A142:C142").Copy()
$ws.Range("A143:C145").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Set final cell contents for every row/column of the table (rows 1-145) ---
$ws.Range("A1").Value = "AkNr"
$ws.Range("B1").Value = "KlartextDe"
$ws.Range("C1").Value = "KlartextEn"
$ws.Range("A2").Value = "z. B."
$ws.Range("B2").Value = "zum Beispiel"
$ws.Range("C2").Value = ""
$ws.Range("A3").Value = "WOB"
$ws.Range("B3").Value = "Women on Board"
$ws.Range("C3").Value = "Women on Board"
$ws.Range("A4").Value = "WHO"
$ws.Range("B4").Value = "Weltgesundheitsorganisation (World Health Organization)"
$ws.Range("C4").Value = "World Health Organization"
$ws.Range("A5").Value = "VN"
$ws.Range("B5").Value = "Vereinte Nationen (United Nations)"
$ws.Range("C5").Value = "United Nations"
$ws.Range("A6").Value = "VGR"
$ws.Range("B6").Value = "Volkswirtschaftlichen Gesamtrechnungen"
$ws.Range("C6").Value = ""
$ws.Range("A7").Value = "v. a."
$ws.Range("B7").Value = "vor allem"
$ws.Range("C7").Value = ""
$ws.Range("A8").Value = "usw."
$ws.Range("B8").Value = "und so weiter"
$ws.Range("C8").Value = ""
$ws.Range("A9").Value = "USD"
$ws.Range("B9").Value = "US-Dollar"
$ws.Range("C9").Value = "United States dollar"
$ws.Range("A10").Value = "USA"
$ws.Range("B10").Value = "Vereinigte Staaten von Amerika (United States of America)"
$ws.Range("C10").Value = "United States of America"
$ws.Range("A11").Value = "US"
$ws.Range("B11").Value = "Vereinigte Staaten von Amerika (United States)"
$ws.Range("C11").Value = "United States"
$ws.Range("A12").Value = "UNFCCC"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "United Nations Framework Convention on Climate Change"
$ws.Range("A13").Value = "UNCCD"
$ws.Range("B13").Value = "Übereinkommens der Vereinten Nationen zur Bekämpfung der Wüstenbildung in Entwicklungs- und Schwellenländern"
$ws.Range("C13").Value = "UN Convention to Combat Desertification"
$ws.Range("A14").Value = "UN"
$ws.Range("B14").Value = "Vereinte Nationen (United Nations)"
$ws.Range("C14").Value = "United Nations"
$ws.Range("A15").Value = "UK"
$ws.Range("B15").Value = "Vereinigtes Königreich (United Kingdom)"
$ws.Range("C15").Value = "United Kingdom"
$ws.Range("A16").Value = "UBA"
$ws.Range("B16").Value = "Umweltbundesamt"
$ws.Range("C16").Value = "Federal Environment Agency"
$ws.Range("A17").Value = "u.a."
$ws.Range("B17").Value = "unter anderem"
$ws.Range("C17").Value = ""
$ws.Range("A18").Value = "u. a."
$ws.Range("B18").Value = "unter anderem"
$ws.Range("C18").Value = ""
$ws.Range("A19").Value = "TWh"
$ws.Range("B19").Value = "Terawattstunde"
$ws.Range("C19").Value = "Terawatt hour"
$ws.Range("A20").Value = "TREMOD"
$ws.Range("B20").Value = "Transport Emission Estimation Model"
$ws.Range("C20").Value = "Transport Emission Estimation Model"
$ws.Range("A21").Value = "TKU"
$ws.Range("B21").Value = "Telekommunikationsunternehmen"
$ws.Range("C21").Value = "Telecommunications companies"
$ws.Range("A22").Value = "SOEP"
$ws.Range("B22").Value = "Sozio-oekonomischen Panel"
$ws.Range("C22").Value = "Socio-Economic Panel"
$ws.Range("A23").Value = "SO₂"
$ws.Range("B23").Value = "Schwefeldioxid"
$ws.Range("C23").Value = "Sulphur dioxide"
$ws.Range("A24").Value = "SMEs"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "Small and medium-sized enterprises"
$ws.Range("A25").Value = "SF₆"
$ws.Range("B25").Value = "Schwefelhexafluorid"
$ws.Range("C25").Value = "Sulphur hexafluoride"
$ws.Range("A26").Value = "SES"
$ws.Range("B26").Value = "Sozioökonomischer Status"
$ws.Range("C26").Value = "Socioeconomic status"
$ws.Range("A27").Value = "SE"
$ws.Range("B27").Value = "Europäische Gesellschaft"
$ws.Range("C27").Value = ""
$ws.Range("A28").Value = "SDGs"
$ws.Range("B28").Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Range("C28").Value = "Sustainable Development Goals"
$ws.Range("A29").Value = "SDG"
$ws.Range("B29").Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Range("C29").Value = "Sustainable Development Goals"
$ws.Range("A30").Value = "SALW"
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = "Small arms and light weapons"
$ws.Range("A31").Value = "RKI"
$ws.Range("B31").Value = "Robert Koch-Institut"
$ws.Range("C31").Value = "Robert Koch-Institute"
$ws.Range("A32").Value = "REDD"
$ws.Range("B32").Value = "Verringerung von Emissionen aus Entwaldung und Waldschädigung sowie die Rolle des Waldschutzes, der nachhaltigen Waldbewirtschaftung und des Ausbaus des Kohlenstoffspeichers Wald in Entwicklungsländern"
$ws.Range("C32").Value = "Reducing emissions from deforestation and forest degradation and the role of conservation, sustainable management of forests and enhancement of forest carbon stocks in developing countries"
$ws.Range("A33").Value = "R&D"
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = "Research and development"
$ws.Range("A34").Value = "PM₂.₅"
$ws.Range("B34").Value = "Feinstaub (Durchmesser kleiner 2,5 Mikrometer)"
$ws.Range("C34").Value = "Particulate matter (diameter smaller than 2.5 micrometers)"
$ws.Range("A35").Value = "PM₂,₅"
$ws.Range("B35").Value = "Feinstaub (Durchmesser kleiner 2,5 Mikrometer)"
$ws.Range("C35").Value = "Particulate matter (diameter smaller than 2.5 micrometers)"
$ws.Range("A36").Value = "PM₁₀"
$ws.Range("B36").Value = "Feinstaub (Durchmesser kleiner 10 Mikrometer)"
$ws.Range("C36").Value = "Particulate matter (diameter smaller than 10 micrometers)"
$ws.Range("A37").Value = "PM₀.₁"
$ws.Range("B37").Value = "Feinstaub (Durchmesser kleiner 0,1 Mikrometer)"
$ws.Range("C37").Value = "Particulate matter (diameter smaller than 0.1 micrometers)"
$ws.Range("A38").Value = "PM₀,₁"
$ws.Range("B38").Value = "Feinstaub (Durchmesser kleiner 0,1 Mikrometer)"
$ws.Range("C38").Value = "Particulate matter (diameter smaller than 0.1 micrometers)"
$ws.Range("A39").Value = "PKS"
$ws.Range("B39").Value = "Polizeilichen Kriminalstatistik"
$ws.Range("C39").Value = "Police Crime Statistics"
$ws.Range("A40").Value = "PINETI"
$ws.Range("B40").Value = "Pollutant INput and EcosysTem Impact"
$ws.Range("C40").Value = "Pollutant INput and EcosysTem Impact"
$ws.Range("A41").Value = "PFCs"
$ws.Range("B41").Value = ""
$ws.Range("C41").Value = "Perfluorocarbons"
$ws.Range("A42").Value = "P97"
$ws.Range("B42").Value = "97. Perzentil"
$ws.Range("C42").Value = "97th percentile"
$ws.Range("A43").Value = "P90"
$ws.Range("B43").Value = "90. Perzentil"
$ws.Range("C43").Value = "90th percentile"
$ws.Range("A44").Value = "OECD"
$ws.Range("B44").Value = "Organisation für wirtschaftliche Zusammenarbeit und Entwicklung (Organisation for Economic Co-operation and Development)"
$ws.Range("C44").Value = "Organisation for Economic Co-operation and Development"
$ws.Range("A45").Value = "ODA"
$ws.Range("B45").Value = "Öffentliche Entwicklungsausgaben (official development assistance)"
$ws.Range("C45").Value = "Official development assistance"
$ws.Range("A46").Value = "Nr."
$ws.Range("B46").Value = "Nummer"
$ws.Range("C46").Value = "Number"
$ws.Range("A47").Value = "NOₓ"
$ws.Range("B47").Value = "Stickstoffoxid"
$ws.Range("C47").Value = "Nitrogen oxides"
$ws.Range("A48").Value = "NMVOCs"
$ws.Range("B48").Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Range("C48").Value = "Non-methane volatile organic compounds"
$ws.Range("A49").Value = "NMVOC"
$ws.Range("B49").Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Range("C49").Value = "non-methane volatile organic compounds"
$ws.Range("A50").Value = "NH₃"
$ws.Range("B50").Value = "Ammoniak"
$ws.Range("C50").Value = "Ammonia"
$ws.Range("A51").Value = "NF₃"
$ws.Range("B51").Value = "Stickstofftrifluorid"
$ws.Range("C51").Value = "Nitrogen trifluoride"
$ws.Range("A52").Value = "NEC"
$ws.Range("B52").Value = "Richtlinie über nationale Emissionshöchstmengen für bestimmte Luftschadstoffe (National Emission Ceilings Directive)"
$ws.Range("C52").Value = "National Emission Ceilings Directive"
$ws.Range("A53").Value = "N₂O"
$ws.Range("B53").Value = "Lachgas"
$ws.Range("C53").Value = "Nitrous oxide"
$ws.Range("A54").Value = "N"
$ws.Range("B54").Value = ""
$ws.Range("C54").Value = "Nitrogen"
$ws.Range("A55").Value = "MSY"
$ws.Range("B55").Value = "Maximum Sustainable Yield"
$ws.Range("C55").Value = "Maximum Sustainable Yield"
$ws.Range("A56").Value = "Mrd."
$ws.Range("B56").Value = "Milliarde"
$ws.Range("C56").Value = ""
$ws.Range("A57").Value = "mg/l"
$ws.Range("B57").Value = "Milligramm pro Liter"
$ws.Range("C57").Value = "Miligrams per litre"
$ws.Range("A58").Value = "mg"
$ws.Range("B58").Value = "Milligramm"
$ws.Range("C58").Value = "Miligrams"
$ws.Range("A59").Value = "Mbps"
$ws.Range("B59").Value = ""
$ws.Range("C59").Value = "Megabit per second"
$ws.Range("A60").Value = "Mbit/s"
$ws.Range("B60").Value = "Megabit pro Sekunde"
$ws.Range("C60").Value = "Megabit per second"
$ws.Range("A61").Value = "m³"
$ws.Range("B61").Value = "Kubikmeter"
$ws.Range("C61").Value = "Cubic metre"
$ws.Range("A62").Value = "m²"
$ws.Range("B62").Value = "Quadratmeter"
$ws.Range("C62").Value = "Square meter"
$ws.Range("A63").Value = "LULUCF"
$ws.Range("B63").Value = ""
$ws.Range("C63").Value = "Land use, land-use change and forestry"
$ws.Range("A64").Value = "LDCs"
$ws.Range("B64").Value = "am wenigsten entwickelte Länder (Least developed countries)"
$ws.Range("C64").Value = "Least developed countries"
$ws.Range("A65").Value = "LDC"
$ws.Range("B65").Value = "am wenigsten entwickelte Länder (Least developed countries)"
$ws.Range("C65").Value = "Least developed countries"
$ws.Range("A66").Value = "LAWA"
$ws.Range("B66").Value = "Bund/Länder-Arbeitsgemeinschaft Wasser"
$ws.Range("C66").Value = "German Working Group on Water Issues of the Länder and the Federal Government"
$ws.Range("A67").Value = "l"
$ws.Range("B67").Value = "Liter"
$ws.Range("C67").Value = "Litre"
$ws.Range("A68").Value = "km²"
$ws.Range("B68").Value = "Quadratkilometer"
$ws.Range("C68").Value = "Square kilometer"
$ws.Range("A69").Value = "KiGGS"
$ws.Range("B69").Value = "Studie zur Gesundheit von Kindern und Jugendlichen in Deutschland"
$ws.Range("C69").Value = "Study on the health of children and adolescents in Germany"
$ws.Range("A70").Value = "kg/m²"
$ws.Range("B70").Value = "Kilogramm pro Quadratmeter"
$ws.Range("C70").Value = "Kilogram per square meter"
$ws.Range("A71").Value = "kg/ha"
$ws.Range("B71").Value = "Kilogramm pro Hektar"
$ws.Range("C71").Value = "Kilogram per hectare"
$ws.Range("A72").Value = "kg"
$ws.Range("B72").Value = "Kilogramm"
$ws.Range("C72").Value = "Kilogram"
$ws.Range("A73").Value = "Kfz"
$ws.Range("B73").Value = "Kraftfahrzeug"
$ws.Range("C73").Value = ""
$ws.Range("A74").Value = "KfW"
$ws.Range("B74").Value = "Kreditanstalt für Wiederaufbau"
$ws.Range("C74").Value = "Reconstruction Loan Corporation (Kreditanstalt für Wiederaufbau)"
$ws.Range("A75").Value = "ISCO"
$ws.Range("B75").Value = "Internationale Standardklassifikation der Berufe (International Standard Classification of Occupations)"
$ws.Range("C75").Value = "International Standard Classification of Occupations"
$ws.Range("A76").Value = "ISCED"
$ws.Range("B76").Value = "Internationale Standardklassifikation des Bildungswesens (International Standard Classification of Education)"
$ws.Range("C76").Value = "International Standard Classification of Education"
$ws.Range("A77").Value = "i.e."
$ws.Range("B77").Value = ""
$ws.Range("C77").Value = "that is to say (id est)"
$ws.Range("A78").Value = "H-FKW/HFC"
$ws.Range("B78").Value = "Teilhalogenierte Fluorkohlenwasserstoffe"
$ws.Range("C78").Value = ""
$ws.Range("A79").Value = "HFCs"
$ws.Range("B79").Value = ""
$ws.Range("C79").Value = "Hydrofluorocarbons"
$ws.Range("A80").Value = "ha"
$ws.Range("B80").Value = "Hektar"
$ws.Range("C80").Value = "Hectare"
$ws.Range("A81").Value = "GPG"
$ws.Range("B81").Value = "Geschlechtsspezifischen Verdienstabstand (gender pay gap)"
$ws.Range("C81").Value = "gender pay gap"
$ws.Range("A82").Value = "GNI"
$ws.Range("B82").Value = ""
$ws.Range("C82").Value = "Gross national income"
$ws.Range("A83").Value = "GmbH"
$ws.Range("B83").Value = "Gesellschaft mit beschränkter Haftung"
$ws.Range("C83").Value = "Company with limited liability"
$ws.Range("A84").Value = "GIZ"
$ws.Range("B84").Value = "Deutsche Gesellschaft für Internationale Zusammenarbeit"
$ws.Range("C84").Value = "German Agency for International Cooperation"
$ws.Range("A85").Value = "gGmbH"
$ws.Range("B85").Value = "gemeinnützige Gesellschaft mit beschränkter Haftung"
$ws.Range("C85").Value = "Non-profit limited liability company"
$ws.Range("A86").Value = "GG"
$ws.Range("B86").Value = "Grundgesetz"
$ws.Range("C86").Value = "Basic Law"
$ws.Range("A87").Value = "GDP"
$ws.Range("B87").Value = ""
$ws.Range("C87").Value = "Gross domestic product"
$ws.Range("A88").Value = "FuE"
$ws.Range("B88").Value = "Forschung und Entwicklung"
$ws.Range("C88").Value = ""
$ws.Range("A89").Value = "FTTB/H"
$ws.Range("B89").Value = "Reine Glasfasernetze"
$ws.Range("C89").Value = "Fully fibre-optic networks"
$ws.Range("A90").Value = "FKW/PFC"
$ws.Range("B90").Value = "Perfluorierte Kohlenwasserstoffe"
$ws.Range("C90").Value = ""
$ws.Range("A91").Value = "FidAR"
$ws.Range("B91").Value = "Frauen in die Aufsichtsräte"
$ws.Range("C91").Value = ""
$ws.Range("A92").Value = "FCPF"
$ws.Range("B92").Value = "Forest Carbon Partnership Facility"
$ws.Range("C92").Value = "Forest Carbon Partnership Facility"
$ws.Range("A93").Value = "FAO"
$ws.Range("B93").Value = "Ernährungs- und Landwirtschaftsorganisation der Vereinten Nationen (Food and Agriculture Organization)"
$ws.Range("C93").Value = "Food and Agriculture Organization"
$ws.Range("A94").Value = "EU-SILC"
$ws.Range("B94").Value = "Statistik über Einkommen und Lebensbedingungen (Statistics on Income and Living Conditions)"
$ws.Range("C94").Value = "Statistics on Income and Living Conditions"
$ws.Range("A95").Value = "EUR"
$ws.Range("B95").Value = "Euro"
$ws.Range("C95").Value = "Euro"
$ws.Range("A96").Value = "EU-EVK"
$ws.Range("B96").Value = "EU-Energieverbrauchskennzeichnung"
$ws.Range("C96").Value = ""
$ws.Range("A97").Value = "EUA"
$ws.Range("B97").Value = "Europäische Umweltagentur"
$ws.Range("C97").Value = ""
$ws.Range("A98").Value = "EU-28"
$ws.Range("B98").Value = "Europäische Union mit 28 Mitgliedsstaaten"
$ws.Range("C98").Value = "European Union consisting of 28 member states"
$ws.Range("A99").Value = "EU-27"
$ws.Range("B99").Value = "Europäische Union mit 27 Mitgliedsstaaten"
$ws.Range("C99").Value = "European Union consisting of 27 member states"
$ws.Range("A100").Value = "EU"
$ws.Range("B100").Value = "Europäische Union"
$ws.Range("C100").Value = "European Union"
$ws.Range("A101").Value = "etc."
$ws.Range("B101").Value = "und so weiter (et cetera)"
$ws.Range("C101").Value = "and so on (et cetera)"
$ws.Range("A102").Value = "ESVG"
$ws.Range("B102").Value = "Europäische System Volkswirtschaftlicher Gesamtrechnungen"
$ws.Range("C102").Value = ""
$ws.Range("A103").Value = "ESA"
$ws.Range("B103").Value = ""
$ws.Range("C103").Value = "European System of National and Regional Accounts"
$ws.Range("A104").Value = "EMAS"
$ws.Range("B104").Value = "Eco-Management and Audit Scheme"
$ws.Range("C104").Value = "Eco-Management and Audit Scheme"
$ws.Range("A105").Value = "EGW"
$ws.Range("B105").Value = "Ernährungs- und der Gewerbliche Wirtschaft"
$ws.Range("C105").Value = "Food and industrial economy"
$ws.Range("A106").Value = "EEG"
$ws.Range("B106").Value = "Erneuerbare-Energien-Gesetz"
$ws.Range("C106").Value = "Renewable Energy Sources Act"
$ws.Range("A107").Value = "EEA"
$ws.Range("B107").Value = ""
$ws.Range("C107").Value = "European Environment Agency"
$ws.Range("A108").Value = "e.g."
$ws.Range("B108").Value = ""
$ws.Range("C108").Value = "for example (exempli gratia)"
$ws.Range("A109").Value = "DIN"
$ws.Range("B109").Value = "Deutsches Institut für Normung e.V."
$ws.Range("C109").Value = "German Institute for Standardisation Registered Association"
$ws.Range("A110").Value = "DEG"
$ws.Range("B110").Value = "Deutsche Investitions- und Entwicklungsgesellschaft"
$ws.Range("C110").Value = "German Investment and Development Corporation (Deutsche Investitions- und Entwicklungsgesellschaft)"
$ws.Range("A111").Value = "DDB"
$ws.Range("B111").Value = "Deutsche Digitale Bibliothek"
$ws.Range("C111").Value = "German Digital Library (Deutsche Digitale Bibliothek)"
$ws.Range("A112").Value = "DDA"
$ws.Range("B112").Value = "Dachverband Deutscher Avifaunisten"
$ws.Range("C112").Value = ""
$ws.Range("A113").Value = "DAC"
$ws.Range("B113").Value = "Richtlinien des Entwicklungsausschusses (Development Assistance Committee)"
$ws.Range("C113").Value = "Development Assistance Committee"
$ws.Range("A114").Value = "CPI"
$ws.Range("B114").Value = "Korruptionswahrnehmungsindex (Corruption Perception Index)"
$ws.Range("C114").Value = "Corruption Perception Index"
$ws.Range("A115").Value = "COVID-19"
$ws.Range("B115").Value = "Coronavirus SARS-CoV-2"
$ws.Range("C115").Value = "Coronavirus SARS-CoV-2"
$ws.Range("A116").Value = "CO₂"
$ws.Range("B116").Value = "Kohlenstoffdioxid"
$ws.Range("C116").Value = "Carbon dioxide"
$ws.Range("A117").Value = "CLRTAP"
$ws.Range("B117").Value = "Genfer Luftreinhaltekonvention (Convention on Long-Range Transboundary Air Pollution)"
$ws.Range("C117").Value = "Convention on Long-Range Transboundary Air Pollution"
$ws.Range("A118").Value = "CH₄"
$ws.Range("B118").Value = "Methan"
$ws.Range("C118").Value = "Methane"
$ws.Range("A119").Value = "CATV"
$ws.Range("B119").Value = "Kabelfernsehen"
$ws.Range("C119").Value = "Cable television"
$ws.Range("A120").Value = "bzw."
$ws.Range("B120").Value = "beziehungsweise"
$ws.Range("C120").Value = ""
$ws.Range("A121").Value = "BNE"
$ws.Range("B121").Value = "Bruttonationaleinkommen"
$ws.Range("C121").Value = ""
$ws.Range("A122").Value = "bn"
$ws.Range("B122").Value = ""
$ws.Range("C122").Value = "Billion"
$ws.Range("A123").Value = "BMZ"
$ws.Range("B123").Value = "Bundesministerium für wirtschaftliche Zusammenarbeit und Entwicklung"
$ws.Range("C123").Value = "Federal Ministry for Economic Cooperation and Development"
$ws.Range("A124").Value = "BMVI"
$ws.Range("B124").Value = "Bundesministerium für Verkehr und digitale Infrastruktur"
$ws.Range("C124").Value = "Federal Ministry of Transport and Digital Infrastructure"
$ws.Range("A125").Value = "BMIs"
$ws.Range("B125").Value = "Body Mass Index"
$ws.Range("C125").Value = "Body Mass Index"
$ws.Range("A126").Value = "BMI"
$ws.Range("B126").Value = "Body Mass Index"
$ws.Range("C126").Value = "Body Mass Index"
$ws.Range("A127").Value = "BMEL"
$ws.Range("B127").Value = "Bundesministeriums für Ernährung und Landwirtschaft"
$ws.Range("C127").Value = "Federal Ministry of Food and Agriculture"
$ws.Range("A128").Value = "BLE"
$ws.Range("B128").Value = "Bundesanstalt für Landwirtschaft und Ernährung"
$ws.Range("C128").Value = "Federal Office for Agriculture and Food"
$ws.Range("A129").Value = "BKG"
$ws.Range("B129").Value = "Bundesamt für Kartographie und Geodäsie"
$ws.Range("C129").Value = "Federal Agency for Cartography and Geodesy"
$ws.Range("A130").Value = "BIP"
$ws.Range("B130").Value = "Bruttoinlandsprodukt"
$ws.Range("C130").Value = ""
$ws.Range("A131").Value = "BfN"
$ws.Range("B131").Value = "Bundesamt für Naturschutz"
$ws.Range("C131").Value = ""
$ws.Range("A132").Value = "BEEG"
$ws.Range("B132").Value = "Bundeselterngeld- und Elternzeitgesetz"
$ws.Range("C132").Value = "Federal Parental Allowance and Parental Leave Act"
$ws.Range("A133").Value = "Art."
$ws.Range("B133").Value = "Artikel"
$ws.Range("C133").Value = "Article"
$ws.Range("A134").Value = "ALKIS"
$ws.Range("B134").Value = "Amtlichen Liegenschaftskataster-Informationssystem"
$ws.Range("C134").Value = "Official land register information system"
$ws.Range("A135").Value = "AGEE-Stat"
$ws.Range("B135").Value = "Arbeitsgruppe Erneuerbare Energien-Statistik"
$ws.Range("C135").Value = "Working Group on Renewable Energy Statistics"
$ws.Range("A136").Value = "AGEB"
$ws.Range("B136").Value = "Arbeitsgemeinschaft Energiebilanzen"
$ws.Range("C136").Value = "Energy Balance Association"
$ws.Range("A137").Value = "AGA"
$ws.Range("B137").Value = "Arbeitsgemeinschaft Adipositas im Kinder- und Jugendalter"
$ws.Range("C137").Value = "Childhood Obesity Federation"
$ws.Range("A138").Value = "Abs."
$ws.Range("B138").Value = "Absatz"
$ws.Range("C138").Value = ""
$ws.Range("A139").Value = "a.m."
$ws.Range("B139").Value = ""
$ws.Range("C139").Value = "before noon (ante meridiem)"
$ws.Range("A140").Value = "µm"
$ws.Range("B140").Value = "Mikrometer"
$ws.Range("C140").Value = "Micrometer"
$ws.Range("A141").Value = "µg/m³"
$ws.Range("B141").Value = "Mikrogramm pro Kubikmeter"
$ws.Range("C141").Value = "Micrograms per cubic metre"
$ws.Range("A142").Value = "µg"
$ws.Range("B142").Value = "Mikrogramm"
$ws.Range("C142").Value = "Micrograms"
$ws.Range("A143").Value = "ICES"
$ws.Range("B143").Value = "Internationalen Rat für Meeresforschung (International Council for the Exploration of the Sea)"
$ws.Range("C143").Value = "International Council for the Exploration of the Sea"
$ws.Range("A144").Value = "CRS"
$ws.Range("B144").Value = "Gläubigermeldesystem (Creditor Reporting System)"
$ws.Range("C144").Value = "`nCreditor Reporting System"
$ws.Range("A145").Value = "PhD"
$ws.Range("B145").Value = ""
$ws.Range("C145").Value = "Doctor of Philosophy (philosophiae doctor)"
